# ByCoach.xlsx update — "Add files via upload"
#
# The workbook's "by Coach" sheet tracks, per player (column A), their coach
# code (column B) and a Yes/No "Started" flag (column C). This edit flips the
# Yes/No flag for a handful of players.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Row -> new Column C text ("Yes"/"No") for every row whose flag flipped.
$updates = @{
    8  = "No"
    9  = "Yes"
    45 = "No"
    46 = "Yes"
    50 = "No"
    51 = "No"
    55 = "Yes"
    57 = "Yes"
    68 = "Yes"
    69 = "No"
    71 = "No"
    72 = "Yes"
    75 = "Yes"
    78 = "No"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Leave the cursor where the author last left it before saving.
$ws.Range("C73").Select() | Out-Null
